$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Cells.Item(17, 8).Value = 6136.6924
$ws.Cells.Item(17, 10).Value = 6136.6924
$ws.Cells.Item(17, 12).Value = 18410.0772
$ws.Cells.Item(17, 14).Value = -18746.0772
# Row 98
$ws.Cells.Item(98, 8).Value = 968.4
$ws.Cells.Item(98, 9).Value = 869.875
$ws.Cells.Item(98, 10).Value = 3333
$ws.Cells.Item(98, 11).Value = 869.875
$ws.Cells.Item(98, 12).Value = 3333
$ws.Cells.Item(98, 13).Value = 628.125
$ws.Cells.Item(98, 14).Value = -6329
# Row 122
$ws.Cells.Item(122, 8).Value = 968.4
$ws.Cells.Item(122, 9).Value = 869.875
$ws.Cells.Item(122, 10).Value = 3333
$ws.Cells.Item(122, 11).Value = 2609.625
$ws.Cells.Item(122, 12).Value = 9999
$ws.Cells.Item(122, 13).Value = -159.625
$ws.Cells.Item(122, 14).Value = -14899
# Row 125
$ws.Cells.Item(125, 8).Value = 68324.56
$ws.Cells.Item(125, 10).Value = 132074.75
$ws.Cells.Item(125, 12).Value = 1188672.75
$ws.Cells.Item(125, 14).Value = -1193592.75

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Cells.Item(4, 8).Value = 451.6
$ws.Cells.Item(4, 9).Value = 389
$ws.Cells.Item(4, 10).Value = 702
$ws.Cells.Item(4, 11).Value = 389
$ws.Cells.Item(4, 12).Value = 702
$ws.Cells.Item(4, 13).Value = -273
$ws.Cells.Item(4, 14).Value = -934
# Row 32
$ws.Cells.Item(32, 8).Value = 9726.103999999999
$ws.Cells.Item(32, 9).Value = 6078.086
$ws.Cells.Item(32, 10).Value = 15277.435
$ws.Cells.Item(32, 11).Value = 6078.086
$ws.Cells.Item(32, 12).Value = 15277.435
$ws.Cells.Item(32, 13).Value = -5791.086
$ws.Cells.Item(32, 14).Value = -15851.435
# Row 61
$ws.Cells.Item(61, 8).Value = 4760.5356
$ws.Cells.Item(61, 9).Value = 4195.522
$ws.Cells.Item(61, 10).Value = 7359.6
$ws.Cells.Item(61, 11).Value = 4195.522
$ws.Cells.Item(61, 12).Value = 7359.6
$ws.Cells.Item(61, 13).Value = -3983.522
$ws.Cells.Item(61, 14).Value = -7783.6
# Row 74
$ws.Cells.Item(74, 8).Value = 4572.1816
$ws.Cells.Item(74, 9).Value = 1980.7
$ws.Cells.Item(74, 10).Value = 6731.75
$ws.Cells.Item(74, 11).Value = 1980.7
$ws.Cells.Item(74, 12).Value = 6731.75
$ws.Cells.Item(74, 13).Value = -1106.7
$ws.Cells.Item(74, 14).Value = -8479.75
# Row 77
$ws.Cells.Item(77, 8).Value = 4572.1816
$ws.Cells.Item(77, 9).Value = 1980.7
$ws.Cells.Item(77, 10).Value = 6731.75
$ws.Cells.Item(77, 11).Value = 9903.5
$ws.Cells.Item(77, 12).Value = 33658.75
$ws.Cells.Item(77, 13).Value = -5535.5
$ws.Cells.Item(77, 14).Value = -42394.75
# Row 97
$ws.Cells.Item(97, 8).Value = 1912.85
$ws.Cells.Item(97, 9).Value = 1937.875
$ws.Cells.Item(97, 10).Value = 1812.75
$ws.Cells.Item(97, 11).Value = 1937.875
$ws.Cells.Item(97, 12).Value = 1812.75
$ws.Cells.Item(97, 13).Value = -1441.875
$ws.Cells.Item(97, 14).Value = -2804.75
# Row 102
$ws.Cells.Item(102, 8).Value = 1990.1
$ws.Cells.Item(102, 9).Value = 1557.4286
$ws.Cells.Item(102, 10).Value = 2999.6667
$ws.Cells.Item(102, 11).Value = 1557.4286
$ws.Cells.Item(102, 12).Value = 2999.6667
$ws.Cells.Item(102, 13).Value = 64.57140000000004
$ws.Cells.Item(102, 14).Value = -6243.6667
# Row 136
$ws.Cells.Item(136, 8).Value = 4760.5356
$ws.Cells.Item(136, 9).Value = 4195.522
$ws.Cells.Item(136, 10).Value = 7359.6
$ws.Cells.Item(136, 11).Value = 12586.566
$ws.Cells.Item(136, 12).Value = 22078.8
$ws.Cells.Item(136, 13).Value = -10036.566
$ws.Cells.Item(136, 14).Value = -27178.8

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Cells.Item(20, 8).Value = 3455.75
$ws.Cells.Item(20, 9).Value = 1712.8
$ws.Cells.Item(20, 10).Value = 4700.7144
$ws.Cells.Item(20, 11).Value = 1712.8
$ws.Cells.Item(20, 12).Value = 4700.7144
$ws.Cells.Item(20, 13).Value = -1465.8
$ws.Cells.Item(20, 14).Value = -5194.7144
# Row 22
$ws.Cells.Item(22, 8).Value = 972.7727
$ws.Cells.Item(22, 9).Value = 850.0625
$ws.Cells.Item(22, 11).Value = 850.0625
$ws.Cells.Item(22, 13).Value = -677.0625
# Row 86
$ws.Cells.Item(86, 8).Value = 2893.3462
$ws.Cells.Item(86, 9).Value = 1183.174
$ws.Cells.Item(86, 11).Value = 1183.174
$ws.Cells.Item(86, 13).Value = -60.17399999999998
# Row 89
$ws.Cells.Item(89, 8).Value = 2893.3462
$ws.Cells.Item(89, 9).Value = 1183.174
$ws.Cells.Item(89, 11).Value = 5915.87
$ws.Cells.Item(89, 13).Value = -299.8699999999999
# Row 105
$ws.Cells.Item(105, 8).Value = 3806.1562
$ws.Cells.Item(105, 10).Value = 3012.2856
$ws.Cells.Item(105, 12).Value = 3012.2856
$ws.Cells.Item(105, 14).Value = -6506.2856
# Row 128
$ws.Cells.Item(128, 8).Value = 36496.668
$ws.Cells.Item(128, 9).Value = 36496.668
$ws.Cells.Item(128, 11).Value = 109490.004
$ws.Cells.Item(128, 13).Value = -107000.004

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Cells.Item(6, 8).Value = 75000
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(6, 11).Value = 0
$ws.Cells.Item(6, 13).Value = ""
# Row 7
$ws.Cells.Item(7, 8).Value = 814.96295
$ws.Cells.Item(7, 9).Value = 793.1111
$ws.Cells.Item(7, 11).Value = 793.1111
$ws.Cells.Item(7, 13).Value = -680.1111
# Row 105
$ws.Cells.Item(105, 8).Value = 2285.25
$ws.Cells.Item(105, 9).Value = 2432.4
$ws.Cells.Item(105, 10).Value = 1549.5
$ws.Cells.Item(105, 11).Value = 2432.4
$ws.Cells.Item(105, 12).Value = 1549.5
$ws.Cells.Item(105, 13).Value = -685.4000000000001
$ws.Cells.Item(105, 14).Value = -5043.5
# Row 134
$ws.Cells.Item(134, 8).Value = 3005.8845
$ws.Cells.Item(134, 9).Value = 2679.7058
$ws.Cells.Item(134, 10).Value = 3622
$ws.Cells.Item(134, 11).Value = 8039.117400000001
$ws.Cells.Item(134, 12).Value = 10866
$ws.Cells.Item(134, 13).Value = -5504.117400000001
$ws.Cells.Item(134, 14).Value = -15936

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = ""
$ws.Cells.Item(7, 14).Value = ""
# Row 8
$ws.Cells.Item(8, 8).Value = 611.06665
$ws.Cells.Item(8, 9).Value = 611.06665
$ws.Cells.Item(8, 11).Value = 1833.19995
$ws.Cells.Item(8, 13).Value = -1694.19995
# Row 11
$ws.Cells.Item(11, 8).Value = 6258506
$ws.Cells.Item(11, 9).Value = 12515629
$ws.Cells.Item(11, 11).Value = 37546887
$ws.Cells.Item(11, 13).Value = -37546747
# Row 43
$ws.Cells.Item(43, 8).Value = 0
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 12).Value = 0
$ws.Cells.Item(43, 14).Value = ""
# Row 51
$ws.Cells.Item(51, 8).Value = 1044.5454
$ws.Cells.Item(51, 10).Value = 999
$ws.Cells.Item(51, 12).Value = 2997
$ws.Cells.Item(51, 14).Value = -3917
# Row 75
$ws.Cells.Item(75, 8).Value = 4016.3
$ws.Cells.Item(75, 9).Value = 2099.5
$ws.Cells.Item(75, 10).Value = 5294.1665
$ws.Cells.Item(75, 11).Value = 6298.5
$ws.Cells.Item(75, 12).Value = 15882.4995
$ws.Cells.Item(75, 13).Value = -5300.5
$ws.Cells.Item(75, 14).Value = -17878.4995
# Row 78
$ws.Cells.Item(78, 8).Value = 4016.3
$ws.Cells.Item(78, 9).Value = 2099.5
$ws.Cells.Item(78, 10).Value = 5294.1665
$ws.Cells.Item(78, 11).Value = 18895.5
$ws.Cells.Item(78, 12).Value = 47647.4985
$ws.Cells.Item(78, 13).Value = -13903.5
$ws.Cells.Item(78, 14).Value = -57631.4985
# Row 92
$ws.Cells.Item(92, 8).Value = 405.13635
$ws.Cells.Item(92, 10).Value = 440.93332
$ws.Cells.Item(92, 12).Value = 1322.79996
$ws.Cells.Item(92, 14).Value = -3818.79996
# Row 103
$ws.Cells.Item(103, 8).Value = 390
$ws.Cells.Item(103, 9).Value = 390
$ws.Cells.Item(103, 10).Value = 0
$ws.Cells.Item(103, 11).Value = 1170
$ws.Cells.Item(103, 12).Value = 0
$ws.Cells.Item(103, 13).Value = -291
$ws.Cells.Item(103, 14).Value = ""
# Row 109
$ws.Cells.Item(109, 8).Value = 1324
$ws.Cells.Item(109, 9).Value = 1324
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 11).Value = 3972
$ws.Cells.Item(109, 12).Value = 0
$ws.Cells.Item(109, 13).Value = -2932
$ws.Cells.Item(109, 14).Value = ""
# Row 118
$ws.Cells.Item(118, 8).Value = 4125.4
$ws.Cells.Item(118, 9).Value = 2656.75
$ws.Cells.Item(118, 11).Value = 7970.25
$ws.Cells.Item(118, 13).Value = -6727.25
# Row 131
$ws.Cells.Item(131, 8).Value = 4175.8184
$ws.Cells.Item(131, 10).Value = 7887.8887
$ws.Cells.Item(131, 12).Value = 23663.6661
$ws.Cells.Item(131, 14).Value = -33743.6661
# Row 133
$ws.Cells.Item(133, 8).Value = 11608.333
$ws.Cells.Item(133, 10).Value = 13166.667
$ws.Cells.Item(133, 12).Value = 39500.001
$ws.Cells.Item(133, 14).Value = -49620.001

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Cells.Item(70, 8).Value = 4999
$ws.Cells.Item(70, 10).Value = 4999
$ws.Cells.Item(70, 12).Value = 4999
$ws.Cells.Item(70, 14).Value = -5539
# Row 73
$ws.Cells.Item(73, 8).Value = 4999
$ws.Cells.Item(73, 10).Value = 4999
$ws.Cells.Item(73, 12).Value = 4999
$ws.Cells.Item(73, 14).Value = -6871
# Row 102
$ws.Cells.Item(102, 8).Value = 1862.4054
$ws.Cells.Item(102, 9).Value = 1211.0714
$ws.Cells.Item(102, 10).Value = 3888.7778
$ws.Cells.Item(102, 11).Value = 1211.0714
$ws.Cells.Item(102, 12).Value = 3888.7778
$ws.Cells.Item(102, 13).Value = 410.9286
$ws.Cells.Item(102, 14).Value = -7132.7778
# Row 132
$ws.Cells.Item(132, 8).Value = 5558.8
$ws.Cells.Item(132, 9).Value = 3978.4736
$ws.Cells.Item(132, 10).Value = 7435.4375
$ws.Cells.Item(132, 11).Value = 11935.4208
$ws.Cells.Item(132, 12).Value = 22306.3125
$ws.Cells.Item(132, 13).Value = -9405.4208
$ws.Cells.Item(132, 14).Value = -27366.3125

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Cells.Item(16, 8).Value = 1219.1666
$ws.Cells.Item(16, 9).Value = 1109.5333
$ws.Cells.Item(16, 10).Value = 1767.3334
$ws.Cells.Item(16, 11).Value = 1109.5333
$ws.Cells.Item(16, 12).Value = 1767.3334
$ws.Cells.Item(16, 13).Value = -939.5333000000001
$ws.Cells.Item(16, 14).Value = -2107.3334
# Row 46
$ws.Cells.Item(46, 8).Value = 3051.2222
$ws.Cells.Item(46, 10).Value = 3638.3572
$ws.Cells.Item(46, 12).Value = 3638.3572
$ws.Cells.Item(46, 14).Value = -4014.3572
# Row 100
$ws.Cells.Item(100, 8).Value = 4930.077
$ws.Cells.Item(100, 9).Value = 4228.7
$ws.Cells.Item(100, 10).Value = 7268
$ws.Cells.Item(100, 11).Value = 4228.7
$ws.Cells.Item(100, 12).Value = 7268
$ws.Cells.Item(100, 13).Value = -3687.7
$ws.Cells.Item(100, 14).Value = -8350
# Row 106
$ws.Cells.Item(106, 8).Value = 111075
$ws.Cells.Item(106, 10).Value = 111075
$ws.Cells.Item(106, 12).Value = 111075
$ws.Cells.Item(106, 14).Value = -113599
# Row 122
$ws.Cells.Item(122, 8).Value = 4788.263
$ws.Cells.Item(122, 9).Value = 4136.933
$ws.Cells.Item(122, 10).Value = 7230.75
$ws.Cells.Item(122, 11).Value = 12410.799
$ws.Cells.Item(122, 12).Value = 21692.25
$ws.Cells.Item(122, 13).Value = -9960.798999999999
$ws.Cells.Item(122, 14).Value = -26592.25
# Row 136
$ws.Cells.Item(136, 8).Value = 6869.5713
$ws.Cells.Item(136, 9).Value = 5546.7144
$ws.Cells.Item(136, 11).Value = 16640.1432
$ws.Cells.Item(136, 13).Value = -14090.1432

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Cells.Item(81, 8).Value = 2629
$ws.Cells.Item(81, 9).Value = 4999
$ws.Cells.Item(81, 10).Value = 1444
$ws.Cells.Item(81, 11).Value = 9998
$ws.Cells.Item(81, 12).Value = 2888
$ws.Cells.Item(81, 13).Value = -8937
$ws.Cells.Item(81, 14).Value = -5010
# Row 84
$ws.Cells.Item(84, 8).Value = 2629
$ws.Cells.Item(84, 9).Value = 4999
$ws.Cells.Item(84, 10).Value = 1444
$ws.Cells.Item(84, 11).Value = 49990
$ws.Cells.Item(84, 12).Value = 14440
$ws.Cells.Item(84, 13).Value = -44686
$ws.Cells.Item(84, 14).Value = -25048
# Row 122
$ws.Cells.Item(122, 8).Value = 2797.2173
$ws.Cells.Item(122, 9).Value = 2301.762
$ws.Cells.Item(122, 11).Value = 6905.286
$ws.Cells.Item(122, 13).Value = -4455.286
